$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Output $ws.Name
Write-Output $ws.Range("A1").Value
